# Update to CCtests Sheets
# Updated to support the new column inversion
#
# Columns I (COORD_X SNAPSHOT GIS (LAT)) and J (COORD_Y SNAPSHOT GIS (LNG))
# were inverted for data rows 2-4, so swap their values back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 4; $r++) {
    $iCell = $ws.Cells.Item($r, 9)
    $jCell = $ws.Cells.Item($r, 10)

    $iVal = $iCell.Value()
    $jVal = $jCell.Value()

    $iCell.Value = $jVal
    $jCell.Value = $iVal
}

# Update the active selection to reflect the last edited cell
$ws.Range("J4").Select()
